# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F4"  = 29
    "F5"  = 429
    "F6"  = 1442
    "F7"  = 944
    "F9"  = 2057
    "F10" = 33
    "F11" = 1239
    "F13" = 95
    "F15" = 299
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
